$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set the Fitness column (C2:C252) values to 7293
$ws.Range("C2:C252").Value = 7293
